# Update "想去人数" (F column) counts across the four worksheets of
# the Beijing comic-convention info workbook, matching the refreshed
# scrape output committed at 456a3b4.

$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions)
$ws = $wb.Worksheets.Item("展览")
$ws.Cells.Item(7, 6).Value  = 5637
$ws.Cells.Item(9, 6).Value  = 7608
$ws.Cells.Item(10, 6).Value = 200
$ws.Cells.Item(11, 6).Value = 67
$ws.Cells.Item(13, 6).Value = 3838
$ws.Cells.Item(16, 6).Value = 198
$ws.Cells.Item(19, 6).Value = 105
$ws.Cells.Item(21, 6).Value = 601
$ws.Cells.Item(22, 6).Value = 3879
$ws.Cells.Item(23, 6).Value = 133
$ws.Cells.Item(25, 6).Value = 5286
$ws.Cells.Item(27, 6).Value = 2092
$ws.Cells.Item(29, 6).Value = 349
$ws.Cells.Item(30, 6).Value = 7854
$ws.Cells.Item(33, 6).Value = 2194
$ws.Cells.Item(34, 6).Value = 2180
$ws.Cells.Item(36, 6).Value = 1293
$ws.Cells.Item(38, 6).Value = 20
$ws.Cells.Item(39, 6).Value = 266
$ws.Cells.Item(40, 6).Value = 247
$ws.Cells.Item(44, 6).Value = 32
$ws.Cells.Item(46, 6).Value = 2070
$ws.Cells.Item(47, 6).Value = 127

# Sheet "演出" (Performances)
$ws = $wb.Worksheets.Item("演出")
$ws.Cells.Item(11, 6).Value = 123
$ws.Cells.Item(20, 6).Value = 7

# Sheet "本地生活" (Local life)
$ws = $wb.Worksheets.Item("本地生活")
$ws.Cells.Item(2, 6).Value = 569
$ws.Cells.Item(3, 6).Value = 746

# Sheet "全部类型" (All types)
$ws = $wb.Worksheets.Item("全部类型")
$ws.Cells.Item(5, 6).Value  = 569
$ws.Cells.Item(6, 6).Value  = 746
$ws.Cells.Item(8, 6).Value  = 5637
$ws.Cells.Item(9, 6).Value  = 7608
$ws.Cells.Item(10, 6).Value = 200
$ws.Cells.Item(11, 6).Value = 3838
$ws.Cells.Item(14, 6).Value = 198
$ws.Cells.Item(17, 6).Value = 105
$ws.Cells.Item(20, 6).Value = 601
$ws.Cells.Item(21, 6).Value = 3879
$ws.Cells.Item(23, 6).Value = 133
$ws.Cells.Item(25, 6).Value = 5286
$ws.Cells.Item(27, 6).Value = 2092
$ws.Cells.Item(29, 6).Value = 349
$ws.Cells.Item(30, 6).Value = 7855
$ws.Cells.Item(33, 6).Value = 2194
$ws.Cells.Item(34, 6).Value = 2180
$ws.Cells.Item(36, 6).Value = 1293
$ws.Cells.Item(37, 6).Value = 266
$ws.Cells.Item(38, 6).Value = 247
$ws.Cells.Item(42, 6).Value = 32
$ws.Cells.Item(44, 6).Value = 2070
$ws.Cells.Item(45, 6).Value = 127
$ws.Cells.Item(48, 6).Value = 7
